$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "MuSCs" sending cluster for this edge is renamed to "Inflammatory-Mac"
$ws.Range("A2").Value = "Inflammatory-Mac"

# Updated TPM-derived numeric values on row 2
$ws.Range("G2").Value = 0.03141166666666666
$ws.Range("H2").Value = 0.094235
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.00006168204277777778
$ws.Range("R2").Value = 0.000555138385
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# Row 3 (duplicate MuSCs->MuSCs edge) is removed entirely, shrinking the
# used range down to A1:T2
$ws.Rows.Item(3).Delete()
